{"js": "// Map of old text -> new text for this edit.\nconst replacements = [\n  [\"2024-12-21 Saturday\", \"2024-12-22 Sunday\"],\n  [\"422\u00d78=\", \"497\u00d74=\"],\n  [\"848\u00d76=\", \"678\u00d74=\"],\n  [\"536\u00d78=\", \"741\u00d76=\"],\n  [\"777\u00d73=\", \"289\u00d77=\"],\n  [\"989\u00d78=\", \"729\u00d76=\"],\n  [\"978\u00d76=\", \"366\u00d77=\"],\n  [\"161\u00d78=\", \"535\u00d75=\"],\n  [\"639\u00d74=\", \"929\u00d73=\"],\n  [\"105\u00d78=\", \"734\u00d73=\"],\n  [\"280\u00d78=\", \"716\u00d74=\"],\n  [\"953\u00d73=\", \"294\u00d72=\"],\n  [\"699\u00d78=\", \"960\u00d79=\"],\n  [\"523\u00d79=\", \"169\u00d79=\"],\n  [\"878\u00d79=\", \"503\u00d79=\"],\n  [\"200\u00d79=\", \"673\u00d72=\"],\n  [\"506\u00d75=\", \"968\u00d78=\"],\n  [\"966\u00d74=\", \"365\u00d76=\"],\n  [\"584\u00d75=\", \"264\u00d78=\"],\n  [\"753\u00d79=\", \"617\u00d75=\"],\n  [\"993\u00d74=\", \"824\u00d75=\"],\n  [\"964\u00d73=\", \"659\u00d76=\"],\n  [\"534\u00d75=\", \"154\u00d72=\"],\n  [\"413\u00d74=\", \"903\u00d78=\"],\n  [\"396\u00d74=\", \"397\u00d74=\"],\n  [\"837\u00d79=\", \"809\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-21 Saturday\", \"2024-12-22 Sunday\"),\n    @(\"422\u00d78=\", \"497\u00d74=\"),\n    @(\"848\u00d76=\", \"678\u00d74=\"),\n    @(\"536\u00d78=\", \"741\u00d76=\"),\n    @(\"777\u00d73=\", \"289\u00d77=\"),\n    @(\"989\u00d78=\", \"729\u00d76=\"),\n    @(\"978\u00d76=\", \"366\u00d77=\"),\n    @(\"161\u00d78=\", \"535\u00d75=\"),\n    @(\"639\u00d74=\", \"929\u00d73=\"),\n    @(\"105\u00d78=\", \"734\u00d73=\"),\n    @(\"280\u00d78=\", \"716\u00d74=\"),\n    @(\"953\u00d73=\", \"294\u00d72=\"),\n    @(\"699\u00d78=\", \"960\u00d79=\"),\n    @(\"523\u00d79=\", \"169\u00d79=\"),\n    @(\"878\u00d79=\", \"503\u00d79=\"),\n    @(\"200\u00d79=\", \"673\u00d72=\"),\n    @(\"506\u00d75=\", \"968\u00d78=\"),\n    @(\"966\u00d74=\", \"365\u00d76=\"),\n    @(\"584\u00d75=\", \"264\u00d78=\"),\n    @(\"753\u00d79=\", \"617\u00d75=\"),\n    @(\"993\u00d74=\", \"824\u00d75=\"),\n    @(\"964\u00d73=\", \"659\u00d76=\"),\n    @(\"534\u00d75=\", \"154\u00d72=\"),\n    @(\"413\u00d74=\", \"903\u00d78=\"),\n    @(\"396\u00d74=\", \"397\u00d74=\"),\n    @(\"837\u00d79=\", \"809\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)\n}\n"}
